# Ajax en el paginador y combo box
# Insert 3 new transaction rows at the top of the statement (rows 1-3),
# pushing all existing rows (and the trailing blank rows) down by 3.
# The H1:H3 "export" formula is a shared formula anchored to rows 1-3 and
# must stay put (it always describes "the first three rows"), so we must
# NOT use a whole-row insert (that would drag the formula down too) -
# instead we shift the A:G values only, using a plain bulk copy.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nbsp = [char]0x00A0

# ---------------------------------------------------------------------
# 1. Shift existing data (rows 1-59, columns A:G) down by 3 rows.
#    Row 59 was the last (blank) row, so this lands on row 62 - matching
#    the new dimension A1:H62.
# ---------------------------------------------------------------------
$srcVals = $ws.Range("A1:G59").Value2
$ws.Range("A4:G62").Value2 = $srcVals

# New rows 60:62 did not exist before, so they need the same formatting
# (date number format on column A) as the other trailing blank rows.
$ws.Range("A59").Copy()
$ws.Range("A60:A62").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Write the 3 new transactions into rows 1-3 (columns A:G only - the
#    H formula for these rows is already in place and will recompute on
#    its own from the new A:G values).
# ---------------------------------------------------------------------
$newRows = @(
    @(41717, "INTERES A SU FAVOR", "C", "0000950820", "AGENCIA PARA PROCESOS BATCH", ("0.26{0}{0}" -f $nbsp), "4111.29"),
    @(41716, "INTERES A SU FAVOR", "C", "0000950832", "AGENCIA PARA PROCESOS BATCH", ("0.26{0}{0}" -f $nbsp), "4111.03"),
    @(41715, "INTERES A SU FAVOR", "C", "0000950846", "AGENCIA PARA PROCESOS BATCH", ("0.26{0}{0}" -f $nbsp), "4110.77")
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $i + 1
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value2 = $row[0]         # A: fecha
    $ws.Cells.Item($r, 2).Value2 = $row[1]         # B: concepto
    $ws.Cells.Item($r, 3).Value2 = $row[2]         # C: tipo
    $ws.Cells.Item($r, 4).Value2 = $row[3]         # D: documento
    $ws.Cells.Item($r, 5).Value2 = $row[4]         # E: oficina
    $ws.Cells.Item($r, 6).Value2 = $row[5]         # F: monto
    $ws.Cells.Item($r, 7).Value2 = $row[6]         # G: saldo
}

$wb.Application.CalculateFull()
